$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) passive values
$ws.Range("B2").Value = 48.282181340544462
$ws.Range("C2").Value = 53.842668214240142
$ws.Range("D2").Value = 51.396987585438829
$ws.Range("E2").Value = 56.412131397130395

# Row 3 (STR) passive values
$ws.Range("B3").Value = 44.044568128126528
$ws.Range("C3").Value = 47.089273951062118
$ws.Range("D3").Value = 47.21481631296556
$ws.Range("E3").Value = 54.032820673687041

# Update the selected range to match the new narrower selection
$ws.Range("B1:E3").Select()
